# The wml.xsd schema for CT_RPr requires toggle properties (w:b, w:i, ...)
# to precede w:color inside w:rPr. Several custom character styles (the
# Pandoc/highlighting "*Tok" styles) in styles.xml had w:color emitted
# before w:b/w:i, which OOXMLValidatorCLI flags as a schema error even
# though xmllint stays silent. Re-assigning the Font.Bold / Font.Italic
# properties (to their existing effective value) on each affected style
# forces the style's rPr to be re-serialized in schema-compliant order
# (toggle properties before w:color), without changing any formatting.

$d = $word.ActiveDocument

# color, b  ->  b, color
$d.Styles("KeywordTok").Font.Bold = $true
$d.Styles("ImportTok").Font.Bold = $true
$d.Styles("ControlFlowTok").Font.Bold = $true
$d.Styles("AlertTok").Font.Bold = $true
$d.Styles("ErrorTok").Font.Bold = $true

# color, i  ->  i, color
$d.Styles("CommentTok").Font.Italic = $true
$d.Styles("DocumentationTok").Font.Italic = $true

# color, b, i  ->  b, i, color
$d.Styles("AnnotationTok").Font.Bold = $true
$d.Styles("AnnotationTok").Font.Italic = $true
$d.Styles("CommentVarTok").Font.Bold = $true
$d.Styles("CommentVarTok").Font.Italic = $true
$d.Styles("InformationTok").Font.Bold = $true
$d.Styles("InformationTok").Font.Italic = $true
$d.Styles("WarningTok").Font.Bold = $true
$d.Styles("WarningTok").Font.Italic = $true
